# supplier_status and comments added to supplier society import functionality
#
# - Supplier Code for the "Jelly Beans" row (E2) changes from ASD to QWE.
# - New "Status" (AB) and "Comments" (AC) values are populated for both
#   data rows.
# - Columns AB/AC get explicit (custom) widths to fit the new content.
# - The sheet view zoom goes from 100% to 120% and the active selection
#   moves from E4 to E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Jelly Beans"): supplier code correction + new status/comments.
$ws.Range("E2").Value = "QWE"
$ws.Range("AB2").Value = "Tapped"
$ws.Range("AC2").Value = "Manager seems to be friendly"

# Row 3 ("Choco Pie"): new status/comments.
$ws.Range("AB3").Value = "LetterGiven"
$ws.Range("AC3").Value = "Manager is not friendly"

# Give the new Status/Comments columns (AB=28, AC=29) explicit widths,
# matching the target ~14.72 / ~8.79 character widths as closely as the
# 1/6-character rounding of ColumnWidth allows.
$ws.Columns.Item(28).ColumnWidth = 13.8333333333
$ws.Columns.Item(29).ColumnWidth = 8.0

# Zoom in on the sheet and move the active selection to E3.
$excel.ActiveWindow.Zoom = 120
$ws.Range("E3").Select() | Out-Null
